# NATMI ligand-receptor edge table refresh (new TPM input) for Wnt8a-Fzd5.
# Updates columns E:T (ligand/receptor stats + derived specificities + edge
# weights) for every data row (2-19); columns A-D and K-L are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Inflammatory-Mac -> ECs
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.172733
$ws.Range("H2").Value = 3.518199
$ws.Range("I2").Value = 0.02094912533655296
$ws.Range("J2").Value = 0.02094912533655295
$ws.Range("M2").Value = 2.815739333333333
$ws.Range("N2").Value = 8.447217999999999
$ws.Range("O2").Value = 0.07700398964630729
$ws.Range("P2").Value = 0.07700398964630729
$ws.Range("Q2").Value = 3.302110435598
$ws.Range("R2").Value = 29.718993920382
$ws.Range("S2").Value = 0.001613166230515118
$ws.Range("T2").Value = 0.001613166230515117

# Row 3: Inflammatory-Mac -> FAPs
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.172733
$ws.Range("H3").Value = 3.518199
$ws.Range("I3").Value = 0.02094912533655296
$ws.Range("J3").Value = 0.02094912533655295
$ws.Range("O3").Value = 0.1324338085883186
$ws.Range("P3").Value = 0.1324338085883186
$ws.Range("Q3").Value = 5.679070180313
$ws.Range("R3").Value = 51.111631622817
$ws.Range("S3").Value = 0.002774372454913749
$ws.Range("T3").Value = 0.002774372454913749

# Row 4: Inflammatory-Mac -> Inflammatory-Mac
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.172733
$ws.Range("H4").Value = 3.518199
$ws.Range("I4").Value = 0.02094912533655296
$ws.Range("J4").Value = 0.02094912533655295
$ws.Range("M4").Value = 5.537790999999999
$ws.Range("N4").Value = 16.613373
$ws.Range("O4").Value = 0.1514458372546134
$ws.Range("P4").Value = 0.1514458372546134
$ws.Range("Q4").Value = 6.494350252803
$ws.Range("R4").Value = 58.449152275227
$ws.Range("S4").Value = 0.003172657826346098
$ws.Range("T4").Value = 0.003172657826346098

# Row 5: Inflammatory-Mac -> MuSCs
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.172733
$ws.Range("H5").Value = 3.518199
$ws.Range("I5").Value = 0.02094912533655296
$ws.Range("J5").Value = 0.02094912533655295
$ws.Range("M5").Value = 1.188595666666667
$ws.Range("N5").Value = 3.565787
$ws.Range("O5").Value = 0.03250535563648733
$ws.Range("P5").Value = 0.03250535563648733
$ws.Range("Q5").Value = 1.393905361957
$ws.Range("R5").Value = 12.545148257613
$ws.Range("S5").Value = 0.0006809587693380011
$ws.Range("T5").Value = 0.000680958769338001

# Row 6: Inflammatory-Mac -> Neutrophils
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.172733
$ws.Range("H6").Value = 3.518199
$ws.Range("I6").Value = 0.02094912533655296
$ws.Range("J6").Value = 0.02094912533655295
$ws.Range("M6").Value = 18.85109966666667
$ws.Range("N6").Value = 56.553299
$ws.Range("O6").Value = 0.5155341854158992
$ws.Range("P6").Value = 0.5155341854158992
$ws.Range("Q6").Value = 22.107306665389
$ws.Range("R6").Value = 198.965759988501
$ws.Range("S6").Value = 0.0107999902655554
$ws.Range("T6").Value = 0.0107999902655554

# Row 7: Inflammatory-Mac -> Resolving-Mac
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.172733
$ws.Range("H7").Value = 3.518199
$ws.Range("I7").Value = 0.02094912533655296
$ws.Range("J7").Value = 0.02094912533655295
$ws.Range("M7").Value = 3.330328666666666
$ws.Range("N7").Value = 9.990985999999999
$ws.Range("O7").Value = 0.09107682345837424
$ws.Range("P7").Value = 0.09107682345837424
$ws.Range("Q7").Value = 3.905586328246
$ws.Range("R7").Value = 35.150276954214
$ws.Range("S7").Value = 0.001907979789884589
$ws.Range("T7").Value = 0.001907979789884588

# Row 8: Neutrophils -> ECs
$ws.Range("G8").Value = 54.58029933333334
$ws.Range("H8").Value = 163.740898
$ws.Range("I8").Value = 0.974995614211059
$ws.Range("J8").Value = 0.974995614211059
$ws.Range("M8").Value = 2.815739333333333
$ws.Range("N8").Value = 8.447217999999999
$ws.Range("O8").Value = 0.07700398964630729
$ws.Range("P8").Value = 0.07700398964630729
$ws.Range("Q8").Value = 153.6838956579738
$ws.Range("R8").Value = 1383.155060921764
$ws.Range("S8").Value = 0.0750785521819034
$ws.Range("T8").Value = 0.0750785521819034

# Row 9: Neutrophils -> FAPs
$ws.Range("G9").Value = 54.58029933333334
$ws.Range("H9").Value = 163.740898
$ws.Range("I9").Value = 0.974995614211059
$ws.Range("J9").Value = 0.974995614211059
$ws.Range("O9").Value = 0.1324338085883186
$ws.Range("P9").Value = 0.1324338085883186
$ws.Range("Q9").Value = 264.3102482632372
$ws.Range("R9").Value = 2378.792234369134
$ws.Range("S9").Value = 0.1291223825468775
$ws.Range("T9").Value = 0.1291223825468775

# Row 10: Neutrophils -> Inflammatory-Mac
$ws.Range("G10").Value = 54.58029933333334
$ws.Range("H10").Value = 163.740898
$ws.Range("I10").Value = 0.974995614211059
$ws.Range("J10").Value = 0.974995614211059
$ws.Range("M10").Value = 5.537790999999999
$ws.Range("N10").Value = 16.613373
$ws.Range("O10").Value = 0.1514458372546134
$ws.Range("P10").Value = 0.1514458372546134
$ws.Range("Q10").Value = 302.2542904254393
$ws.Range("R10").Value = 2720.288613828954
$ws.Range("S10").Value = 0.1476590271137699
$ws.Range("T10").Value = 0.1476590271137699

# Row 11: Neutrophils -> MuSCs
$ws.Range("G11").Value = 54.58029933333334
$ws.Range("H11").Value = 163.740898
$ws.Range("I11").Value = 0.974995614211059
$ws.Range("J11").Value = 0.974995614211059
$ws.Range("M11").Value = 1.188595666666667
$ws.Range("N11").Value = 3.565787
$ws.Range("O11").Value = 0.03250535563648733
$ws.Range("P11").Value = 0.03250535563648733
$ws.Range("Q11").Value = 64.87390727296956
$ws.Range("R11").Value = 583.8651654567261
$ws.Range("S11").Value = 0.03169257918394587
$ws.Range("T11").Value = 0.03169257918394587

# Row 12: Neutrophils -> Neutrophils
$ws.Range("G12").Value = 54.58029933333334
$ws.Range("H12").Value = 163.740898
$ws.Range("I12").Value = 0.974995614211059
$ws.Range("J12").Value = 0.974995614211059
$ws.Range("M12").Value = 18.85109966666667
$ws.Range("N12").Value = 56.553299
$ws.Range("O12").Value = 0.5155341854158992
$ws.Range("P12").Value = 0.5155341854158992
$ws.Range("Q12").Value = 1028.898662569167
$ws.Range("R12").Value = 9260.087963122503
$ws.Range("S12").Value = 0.5026435697563726
$ws.Range("T12").Value = 0.5026435697563726

# Row 13: Neutrophils -> Resolving-Mac
$ws.Range("G13").Value = 54.58029933333334
$ws.Range("H13").Value = 163.740898
$ws.Range("I13").Value = 0.974995614211059
$ws.Range("J13").Value = 0.974995614211059
$ws.Range("M13").Value = 3.330328666666666
$ws.Range("N13").Value = 9.990985999999999
$ws.Range("O13").Value = 0.09107682345837424
$ws.Range("P13").Value = 0.09107682345837424
$ws.Range("Q13").Value = 181.7703355050476
$ws.Range("R13").Value = 1635.933019545428
$ws.Range("S13").Value = 0.08879950342818978
$ws.Range("T13").Value = 0.08879950342818978

# Row 14: Resolving-Mac -> ECs
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.2270136666666667
$ws.Range("H14").Value = 0.681041
$ws.Range("I14").Value = 0.0040552604523881
$ws.Range("J14").Value = 0.0040552604523881
$ws.Range("M14").Value = 2.815739333333333
$ws.Range("N14").Value = 8.447217999999999
$ws.Range("O14").Value = 0.07700398964630729
$ws.Range("P14").Value = 0.07700398964630729
$ws.Range("Q14").Value = 0.6392113104375555
$ws.Range("R14").Value = 5.752901793937999
$ws.Range("S14").Value = 0.0003122712338887727
$ws.Range("T14").Value = 0.0003122712338887727

# Row 15: Resolving-Mac -> FAPs
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.2270136666666667
$ws.Range("H15").Value = 0.681041
$ws.Range("I15").Value = 0.0040552604523881
$ws.Range("J15").Value = 0.0040552604523881
$ws.Range("O15").Value = 0.1324338085883186
$ws.Range("P15").Value = 0.1324338085883186
$ws.Range("Q15").Value = 1.099335095789222
$ws.Range("R15").Value = 9.894015862103
$ws.Range("S15").Value = 0.0005370535865273438
$ws.Range("T15").Value = 0.0005370535865273438

# Row 16: Resolving-Mac -> Inflammatory-Mac
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.2270136666666667
$ws.Range("H16").Value = 0.681041
$ws.Range("I16").Value = 0.0040552604523881
$ws.Range("J16").Value = 0.0040552604523881
$ws.Range("M16").Value = 5.537790999999999
$ws.Range("N16").Value = 16.613373
$ws.Range("O16").Value = 0.1514458372546134
$ws.Range("P16").Value = 0.1514458372546134
$ws.Range("Q16").Value = 1.257154240143667
$ws.Range("R16").Value = 11.314388161293
$ws.Range("S16").Value = 0.0006141523144974383
$ws.Range("T16").Value = 0.0006141523144974383

# Row 17: Resolving-Mac -> MuSCs
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.2270136666666667
$ws.Range("H17").Value = 0.681041
$ws.Range("I17").Value = 0.0040552604523881
$ws.Range("J17").Value = 0.0040552604523881
$ws.Range("M17").Value = 1.188595666666667
$ws.Range("N17").Value = 3.565787
$ws.Range("O17").Value = 0.03250535563648733
$ws.Range("P17").Value = 0.03250535563648733
$ws.Range("Q17").Value = 0.2698274604741112
$ws.Range("R17").Value = 2.428447144267
$ws.Range("S17").Value = 0.0001318176832034577
$ws.Range("T17").Value = 0.0001318176832034577

# Row 18: Resolving-Mac -> Neutrophils
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 0.3333333333333333
$ws.Range("G18").Value = 0.2270136666666667
$ws.Range("H18").Value = 0.681041
$ws.Range("I18").Value = 0.0040552604523881
$ws.Range("J18").Value = 0.0040552604523881
$ws.Range("M18").Value = 18.85109966666667
$ws.Range("N18").Value = 56.553299
$ws.Range("O18").Value = 0.5155341854158992
$ws.Range("P18").Value = 0.5155341854158992
$ws.Range("Q18").Value = 4.279457256028778
$ws.Range("R18").Value = 38.51511530425901
$ws.Range("S18").Value = 0.00209062539397121
$ws.Range("T18").Value = 0.00209062539397121

# Row 19: Resolving-Mac -> Resolving-Mac
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 0.3333333333333333
$ws.Range("G19").Value = 0.2270136666666667
$ws.Range("H19").Value = 0.681041
$ws.Range("I19").Value = 0.0040552604523881
$ws.Range("J19").Value = 0.0040552604523881
$ws.Range("M19").Value = 3.330328666666666
$ws.Range("N19").Value = 9.990985999999999
$ws.Range("O19").Value = 0.09107682345837424
$ws.Range("P19").Value = 0.09107682345837424
$ws.Range("Q19").Value = 0.7560301218251111
$ws.Range("R19").Value = 6.804271096426
$ws.Range("S19").Value = 0.0003693402402998778
$ws.Range("T19").Value = 0.0003693402402998778

Write-Output "Updated 234 cells across 18 rows"
